# feat: add 2022-Q3 data
#
# The workbook originally has two sheets: "总计" (totals) and "2021-Q1"
# (fund-holder detail for that quarter). This script:
#   1. Inserts a new row into "总计" for the "2022-Q3" quarter (pushing the
#      existing "2021-Q1" row down).
#   2. Replaces the "2021-Q1" detail sheet with two sheets, in order:
#      "2022-Q3" (new fund-holder detail) and "2021-Q1" (the original
#      fund-holder detail, recreated verbatim) so that sheet order becomes
#      总计, 2022-Q3, 2021-Q1.

$xlPasteFormats = -4122

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Step 1: update the "总计" (totals) summary sheet
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item(1)

# Push the existing data row (2021-Q1) from row 2 down to row 3.
$summary.Rows(2).Insert()
# Row-insert in this engine copies some formatting from the row above onto
# the freshly inserted blank row; clear it so B2:D2 have no explicit style,
# matching the plain (unstyled) data cells used elsewhere in this sheet.
$summary.Range("B2:D2").ClearFormats()

# Give the new A2 cell the same style as the index cells in column A
# (bold/centered/bordered - style used by A3, which used to be A2).
$summary.Range("A3").Copy() | Out-Null
$summary.Range("A2").PasteSpecial($xlPasteFormats)

# New row for the 2022-Q3 quarter.
$summary.Range("A2").Value = 0
$summary.Range("B2").Value = "2022-Q3"
$summary.Range("C2").Value = 4
$summary.Range("D2").Value = 0.13

# Renumber the shifted 2021-Q1 row's index cell.
$summary.Range("A3").Value = 1

# ---------------------------------------------------------------------
# Step 2: remove the old "2021-Q1" detail sheet
# ---------------------------------------------------------------------
$old2021 = $wb.Worksheets.Item(2)
$old2021.Delete() | Out-Null

# ---------------------------------------------------------------------
# Step 3: create the new "2022-Q3" detail sheet right after "总计"
# ---------------------------------------------------------------------
$q3 = $wb.Worksheets.Add($null, $summary)
$q3.Name = "2022-Q3"

# Header row (style copied from the summary sheet's index-cell style so it
# matches the bold/centered/bordered header look used across this workbook).
$summary.Range("A3").Copy() | Out-Null
$q3.Range("B1:H1").PasteSpecial($xlPasteFormats)

$q3.Range("B1").Value = "基金代码"
$q3.Range("C1").Value = "基金名称"
$q3.Range("D1").Value = "基金规模"
$q3.Range("E1").Value = "股票总仓位"
$q3.Range("F1").Value = "仓位占比"
$q3.Range("G1").Value = "持有市值(亿元)"
$q3.Range("H1").Value = "仓位排名"

# Column A index cells use the same bold/centered/bordered style.
$summary.Range("A3").Copy() | Out-Null
$q3.Range("A2:A5").PasteSpecial($xlPasteFormats)

# Cells that hold numeric-looking text (fund codes, percentages, etc.) must
# stay text, not be silently reinterpreted as numbers.
$q3.Range("B2:G5").NumberFormat = "@"

$q3Rows = @(
    @(0, "501030", "汇添富中证环境治理指数（LOF）A", "3.12", "92.74", "2.19", "0.0683", 7),
    @(1, "164908", "交银施罗德中证环境治理指数（LOF）", "1.57", "93.62", "2.17", "0.0341", 7),
    @(2, "501031", "汇添富中证环境治理指数（LOF）C", "1.30", "92.74", "2.19", "0.0285", 7),
    @(3, "013413", "交银施罗德中证环境治理指数（LOF）C", "0.09", "93.62", "2.17", "0.0020", 7)
)

$r = 2
foreach ($row in $q3Rows) {
    $q3.Cells.Item($r, 1).Value = $row[0]
    $q3.Cells.Item($r, 2).Value = $row[1]
    $q3.Cells.Item($r, 3).Value = $row[2]
    $q3.Cells.Item($r, 4).Value = $row[3]
    $q3.Cells.Item($r, 5).Value = $row[4]
    $q3.Cells.Item($r, 6).Value = $row[5]
    $q3.Cells.Item($r, 7).Value = $row[6]
    $q3.Cells.Item($r, 8).Value = $row[7]
    $r = $r + 1
}

$q3.Range("A1").Select()

# ---------------------------------------------------------------------
# Step 4: recreate the "2021-Q1" detail sheet right after "2022-Q3"
# ---------------------------------------------------------------------
$q1 = $wb.Worksheets.Add($null, $q3)
$q1.Name = "2021-Q1"

$summary.Range("A3").Copy() | Out-Null
$q1.Range("B1:H1").PasteSpecial($xlPasteFormats)

$q1.Range("B1").Value = "基金代码"
$q1.Range("C1").Value = "基金名称"
$q1.Range("D1").Value = "基金金额"
$q1.Range("E1").Value = "股票总仓位"
$q1.Range("F1").Value = "仓位占比"
$q1.Range("G1").Value = "持有市值(亿元)"
$q1.Range("H1").Value = "仓位排名"

$summary.Range("A3").Copy() | Out-Null
$q1.Range("A2:A3").PasteSpecial($xlPasteFormats)

$q1.Range("B2:G3").NumberFormat = "@"

$q1Rows = @(
    @(0, "001914", "中信建投聚利混合A", "0.20", "39.92", "2.68", "0.0054", 1),
    @(1, "006845", "中信建投聚利混合C", "0.03", "39.92", "2.68", "0.0008", 1)
)

$r = 2
foreach ($row in $q1Rows) {
    $q1.Cells.Item($r, 1).Value = $row[0]
    $q1.Cells.Item($r, 2).Value = $row[1]
    $q1.Cells.Item($r, 3).Value = $row[2]
    $q1.Cells.Item($r, 4).Value = $row[3]
    $q1.Cells.Item($r, 5).Value = $row[4]
    $q1.Cells.Item($r, 6).Value = $row[5]
    $q1.Cells.Item($r, 7).Value = $row[6]
    $q1.Cells.Item($r, 8).Value = $row[7]
    $r = $r + 1
}

$q1.Range("A1").Select()

Write-Host "2022-Q3 sheet added"
